$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$ws = $wb.Worksheets.Item("P_valores")

$ws.Range("C2").Value = 0.7183209646907236
$ws.Range("D2").Value = 0.7402249843600583
$ws.Range("E2").Value = 0.3268422071057444
$ws.Range("F2").Value = 0.6988980758204226

$ws.Range("B3").Value = 0.7183209646907236
$ws.Range("D3").Value = 0.8683246130014679
$ws.Range("E3").Value = 0.6445270345066689
$ws.Range("F3").Value = 0.9405808542500633

$ws.Range("B4").Value = 0.7402249843600583
$ws.Range("C4").Value = 0.8683246130014679
$ws.Range("E4").Value = 0.3729973656811829
$ws.Range("F4").Value = 0.8878648575023229

$ws.Range("B5").Value = 0.3268422071057444
$ws.Range("C5").Value = 0.6445270345066689
$ws.Range("D5").Value = 0.3729973656811829
$ws.Range("F5").Value = 0.5252335178654639

$ws.Range("B6").Value = 0.6988980758204226
$ws.Range("C6").Value = 0.9405808542500633
$ws.Range("D6").Value = 0.8878648575023229
$ws.Range("E6").Value = 0.5252335178654639

# --- Sheet: Estadisticos_DM ---
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")

$ws2.Range("C2").Value = 0.3680775548080836
$ws2.Range("D2").Value = 0.3382076904650259
$ws2.Range("E2").Value = 1.016054171212341
$ws2.Range("F2").Value = 0.3948578203938488

$ws2.Range("B3").Value = -0.3680775548080836
$ws2.Range("D3").Value = -0.1688568304035477
$ws2.Range("E3").Value = 0.4715291880476259
$ws2.Range("F3").Value = -0.07588945762004497

$ws2.Range("B4").Value = -0.3382076904650259
$ws2.Range("C4").Value = 0.1688568304035477
$ws2.Range("E4").Value = 0.9203094369632598
$ws2.Range("F4").Value = 0.1435973663658751

$ws2.Range("B5").Value = -1.016054171212341
$ws2.Range("C5").Value = -0.4715291880476259
$ws2.Range("D5").Value = -0.9203094369632598
$ws2.Range("F5").Value = -0.6515659294572387

$ws2.Range("B6").Value = -0.3948578203938488
$ws2.Range("C6").Value = 0.07588945762004497
$ws2.Range("D6").Value = -0.1435973663658751
$ws2.Range("E6").Value = 0.6515659294572387
